# Bunch of bug fixes and updates to files & documentation
#
# This edit:
#   1. Removes the two example student rows from the "Data" sheet (the
#      "ME-MS"/"Schmoe, Joe" and "Doe, John"/"ME-PhD" sample rows), leaving
#      only the header row.
#   2. Adds a new "Notes" worksheet (placed before "Data") that documents
#      what each column of the Data sheet means.
#   3. Tweaks some column widths / selected cells to match the new layout.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Clean up the "Data" sheet - remove the two sample/example rows.
# ---------------------------------------------------------------------
$dataSheet = $wb.Worksheets.Item("Data")

$null = $dataSheet.Rows.Item(2).Delete()
$null = $dataSheet.Rows.Item(2).Delete()

# Widen the first two columns now that the sheet only has header text.
$dataSheet.Columns.Item(1).ColumnWidth = 16.8
$dataSheet.Columns.Item(2).ColumnWidth = 17.6

# Matches the new saved selection on the Data sheet.
$null = $dataSheet.Range("C7").Select()

# ---------------------------------------------------------------------
# 2. Add the new "Notes" sheet explaining the fields, inserted before
#    the Data sheet (Worksheets.Add() inserts before the active sheet).
# ---------------------------------------------------------------------
$notes = $wb.Worksheets.Add()
$notes.Name = "Notes"

$notes.Range("A1").Value = "Notes"

$notes.Range("A3").Value = "Fields"

$notes.Range("A4").Value = "Student Name"
$notes.Range("B4").Value = "Any format is fine, probably last, first is best"

$notes.Range("A5").Value = "Current Program"
$notes.Range("B5").Value = "Typically MS or PhD"

$notes.Range("A6").Value = "Start Date"
$notes.Range("B6").Value = "Start date in program"

$notes.Columns.Item(1).ColumnWidth = 17.1

$null = $notes.Range("A6").Select()
